$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (A2 stays 0) - zero-out stat columns, clear summoner name, keep SOLO role
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = ""
$ws.Range("H2").Value = "SOLO"
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0

# Update row 3 with new stats for Portgas D Åce
$ws.Range("B3").Value = 2.719359230470012
$ws.Range("C3").Value = 3811.5
$ws.Range("D3").Value = 0.01770923299731227
$ws.Range("E3").Value = 28.5
$ws.Range("F3").Value = 279.5
$ws.Range("G3").Value = "Portgas D Åce "
$ws.Range("H3").Value = "SOLO"
$ws.Range("I3").Value = 0.1933071592924438
$ws.Range("J3").Value = 18
$ws.Range("K3").Value = 0.01137944447740923

# Remove rows 4-8 (old data for LS DUFFY, BigFather Rengar, Booogeyman, Cantare, Poppy Gods)
$ws.Range("A4:K8").Delete()
